$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column C (authored target width is 39.125 chars; this engine quantizes
# ColumnWidth to 1/7-character steps, so 38.43 is the closest reachable input)
$ws.Columns.Item(3).ColumnWidth = 38.43

# Fill in the new Subject / External-Data-Source cells for the 2020 rows (7-12).
# Cells are written in the same left-to-right / top-to-bottom order the
# original author typed them in, so that newly-introduced shared strings land
# at the same table indices as the authored workbook.

# Row 7
$ws.Range("C7").Value = "Trends in Data Professionals"
$ws.Range("D7").Value = "None?"

# Row 8
$ws.Range("D8").Value = "Anaconda Report, Many"
$ws.Range("C8").Value = "Enthusiast to Data Professional"

# Row 9
$ws.Range("C9").Value = "How to become top data ppl"
$ws.Range("D9").Value = "Macdonalds local prices"

# Row 10
$ws.Range("D10").Value = "Stack Overflow survey"
$ws.Range("C10").Value = "Tools preferences (Kaggle & Stack Overflow)"

# Row 11
$ws.Range("C11").Value = "Trends…"
$ws.Range("D11").Value = "None?"

# Row 12
$ws.Range("C12").Value = "Education Level"
$ws.Range("D12").Value = "None?"

# Match the existing table formatting (center / center / wrap) for the newly
# populated cells.
$ws.Range("C7:D12").HorizontalAlignment = -4108
$ws.Range("C7:D12").VerticalAlignment = -4108
$ws.Range("C7:D12").WrapText = $true

# Rows grew taller once the C/D columns had content in them.
$ws.Rows.Item(7).RowHeight = 21
$ws.Rows.Item(8).RowHeight = 21
$ws.Rows.Item(9).RowHeight = 21
$ws.Rows.Item(10).RowHeight = 21
$ws.Rows.Item(11).RowHeight = 21
$ws.Rows.Item(12).RowHeight = 21

# Leave the selection where the author ended up.
$ws.Range("D13").Select()
